# Add a "cfop" breakdown, mirroring the existing "dept" breakdown:
#   1. "PI hours" gets a new column G "cfop" with the per-person cfop list.
#   2. A new "cfop hours" sheet is appended (after "unit(accumulative) hours")
#      summarizing hours/percentage per cfop value, just like the existing
#      "department hours" / "unit(accumulative) hours" sheets do for dept/unit.

$wb = $excel.ActiveWorkbook
$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1. "PI hours": add column G "cfop".
# ---------------------------------------------------------------------
$piSheet = $wb.Worksheets.Item("PI hours")

# Reuse the existing header formatting (bold + border) from column F.
$piSheet.Range("F1").Copy()
$piSheet.Range("G1").PasteSpecial($xlPasteFormats)

$piSheet.Range("G1").Value = "cfop"
$piSheet.Range("G2").Value = "['cfop_NH']"
$piSheet.Range("G3").Value = "['cfop_ANSELL']"

# ---------------------------------------------------------------------
# 2. New "cfop hours" sheet, appended after the last existing sheet.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$cfopSheet = $wb.Worksheets.Add($null, $lastSheet)
$cfopSheet.Name = "cfop hours"

# Reuse the existing header (bold + border) and index-column formatting.
$piSheet.Range("B1").Copy()
$cfopSheet.Range("B1:D1").PasteSpecial($xlPasteFormats)
$piSheet.Range("A2").Copy()
$cfopSheet.Range("A2:A3").PasteSpecial($xlPasteFormats)

$cfopSheet.Range("B1").Value = "cfop"
$cfopSheet.Range("C1").Value = "hours"
$cfopSheet.Range("D1").Value = "percentage"

$cfopSheet.Range("A2").Value = 0
$cfopSheet.Range("B2").Value = "cfop_NH"
$cfopSheet.Range("C2").Value = 8
$cfopSheet.Range("D2").Value = 72.72727272727273

$cfopSheet.Range("A3").Value = 1
$cfopSheet.Range("B3").Value = "cfop_ANSELL"
$cfopSheet.Range("C3").Value = 3
$cfopSheet.Range("D3").Value = 27.27272727272727

# Creating the sheet made it active; restore the original active tab.
$piSheet.Activate()
